$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F4").Value = 1.54
$ws.Range("G4").Value = 1.94
$ws.Range("I4").Value = 16.5
$ws.Range("K4").Value = 6.8
$ws.Range("V4").Value = 1.06
$ws.Range("W4").Value = 2.06
$ws.Range("L6").Value = 1.38
$ws.Range("S9").Value = 2.6
$ws.Range("N10").Value = 3.05
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 3.05
$ws.Range("H11").Value = 2.72
$ws.Range("I11").Value = 2.74
$ws.Range("R11").Value = 1.21
$ws.Range("U11").Value = 1.81
$ws.Range("V11").Value = 1.57
$ws.Range("W11").Value = 1.48
$ws.Range("X11").Value = 9
$ws.Range("M12").Value = 1.08
$ws.Range("L13").Value = 1.38
$ws.Range("L14").Value = 1.44
$ws.Range("AH14").Value = 19.5
$ws.Range("L18").Value = 1.25
$ws.Range("S18").Value = 2.06
$ws.Range("AB18").Value = 26
$ws.Range("AF18").Value = 34
$ws.Range("AG18").Value = 16
$ws.Range("AN19").Value = 8.199999999999999
$ws.Range("H20").Value = 13.5
$ws.Range("I20").Value = 14
$ws.Range("L20").Value = 1.25
$ws.Range("P20").Value = 2.82
$ws.Range("V20").Value = 1.07
$ws.Range("W20").Value = 4.4
$ws.Range("AA20").Value = 490
$ws.Range("N21").Value = 4.9
$ws.Range("K23").Value = 110

Write-Output "Done applying changes"
